# Scheduled runner update: refresh market-price snapshot columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the rows
# whose underlying item prices changed, across all 8 Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4999.8213
$ws.Range("I98").Value = 4999.8213
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4999.8213
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H121").Value = 1439.8334
$ws.Range("J121").Value = 1434.3636
$ws.Range("L121").Value = 4303.0908
$ws.Range("N121").Value = -7797.0908
$ws.Range("H122").Value = 4999.8213
$ws.Range("I122").Value = 4999.8213
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14999.4639
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H129").Value = 816.95557
$ws.Range("I129").Value = 334.2143
$ws.Range("J129").Value = 1034.9678
$ws.Range("K129").Value = 1002.6429
$ws.Range("L129").Value = 3104.9034
$ws.Range("M129").Value = 3997.3571
$ws.Range("N129").Value = -13104.9034
$ws.Range("H132").Value = 6293600.5
$ws.Range("I132").Value = 9525536
$ws.Range("J132").Value = 9281.888999999999
$ws.Range("K132").Value = 28576608
$ws.Range("L132").Value = 27845.667
$ws.Range("M132").Value = -28574078
$ws.Range("N132").Value = -32905.667
$ws.Range("H135").Value = 1005.55554
$ws.Range("I135").Value = 431.5
$ws.Range("K135").Value = 3883.5
$ws.Range("M135").Value = -1348.5
$ws.Range("H138").Value = 1454.52
$ws.Range("I138").Value = 680.4545000000001
$ws.Range("J138").Value = 1672.8462
$ws.Range("K138").Value = 2041.3635
$ws.Range("L138").Value = 5018.5386
$ws.Range("M138").Value = 3098.6365
$ws.Range("N138").Value = -15298.5386
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 595.5
$ws.Range("I141").Value = 543.5909
$ws.Range("J141").Value = 1737.5
$ws.Range("K141").Value = 1630.7727
$ws.Range("L141").Value = 5212.5
$ws.Range("M141").Value = 3549.2273
$ws.Range("N141").Value = -15572.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 3411.0142
$ws.Range("I32").Value = 3219.4329
$ws.Range("K32").Value = 3219.4329
$ws.Range("M32").Value = -2932.4329
$ws.Range("H63").Value = 2471.1538
$ws.Range("I63").Value = 2133.3333
$ws.Range("J63").Value = 2572.5
$ws.Range("K63").Value = 2133.3333
$ws.Range("L63").Value = 2572.5
$ws.Range("M63").Value = -1447.3333
$ws.Range("N63").Value = -3944.5
$ws.Range("H66").Value = 2471.1538
$ws.Range("I66").Value = 2133.3333
$ws.Range("J66").Value = 2572.5
$ws.Range("K66").Value = 10666.6665
$ws.Range("L66").Value = 12862.5
$ws.Range("M66").Value = -7234.666499999999
$ws.Range("N66").Value = -19726.5
$ws.Range("H97").Value = 494.13333
$ws.Range("I97").Value = 393.23077
$ws.Range("K97").Value = 393.23077
$ws.Range("M97").Value = 102.76923
$ws.Range("H132").Value = 4024.2942
$ws.Range("I132").Value = 6131.857
$ws.Range("J132").Value = 2549
$ws.Range("K132").Value = 18395.571
$ws.Range("L132").Value = 7647
$ws.Range("M132").Value = -15865.571
$ws.Range("N132").Value = -12707
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H82").Value = 30222.223
$ws.Range("J82").Value = 30285.715
$ws.Range("L82").Value = 30285.715
$ws.Range("N82").Value = -31051.715
$ws.Range("H85").Value = 30222.223
$ws.Range("J85").Value = 30285.715
$ws.Range("L85").Value = 30285.715
$ws.Range("N85").Value = -32937.715
$ws.Range("H94").Value = 19232600
$ws.Range("I94").Value = 27779242
$ws.Range("K94").Value = 27779242
$ws.Range("M94").Value = -27778791
$ws.Range("H111").Value = 1500
$ws.Range("J111").Value = 1500
$ws.Range("L111").Value = 1500
$ws.Range("N111").Value = -9680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3065.8462
$ws.Range("I31").Value = 3529.2
$ws.Range("K31").Value = 3529.2
$ws.Range("M31").Value = -3234.2
$ws.Range("H34").Value = 3065.8462
$ws.Range("I34").Value = 3529.2
$ws.Range("K34").Value = 3529.2
$ws.Range("M34").Value = -3327.2
$ws.Range("H141").Value = 29560
$ws.Range("J141").Value = 29560
$ws.Range("L141").Value = 29560
$ws.Range("N141").Value = -39920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 89.666664
$ws.Range("I8").Value = 89.666664
$ws.Range("K8").Value = 268.999992
$ws.Range("M8").Value = -129.999992
$ws.Range("J86").Value = 690
$ws.Range("L86").Value = 2070
$ws.Range("N86").Value = -4442
$ws.Range("J89").Value = 690
$ws.Range("L89").Value = 6210
$ws.Range("N89").Value = -18066
$ws.Range("H122").Value = 798.70966
$ws.Range("J122").Value = 822.4167
$ws.Range("L122").Value = 7401.7503
$ws.Range("N122").Value = -12301.7503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 166.91667
$ws.Range("J2").Value = 352.33334
$ws.Range("L2").Value = 352.33334
$ws.Range("N2").Value = -578.33334
$ws.Range("H132").Value = 5424.25
$ws.Range("I132").Value = 5849.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 17548.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -15018.5
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H46").Value = 5238.8887
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 5837.5
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 5837.5
$ws.Range("M46").Value = -262
$ws.Range("N46").Value = -6213.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1284
$ws.Range("H28").Value = 45009.5
$ws.Range("J28").Value = 70019
$ws.Range("L28").Value = 70019
$ws.Range("N28").Value = -70715
